$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds crypto prices as plain text (values like "54.361.54" use
# "." as a thousands separator and are not valid numbers). Some updated
# values (e.g. "22.40") DO look numeric, and Excel would silently coerce
# them to a Number (dropping the significant trailing zero) unless the cell
# is explicitly formatted as Text first. The NumberFormat is restored to the
# default "Normal" style afterwards so no stray formatting is left behind.

$ws.Range("D2").Value = '54.361.54'
$ws.Range("E2").Value = '  +1.41%  '
$ws.Range("D3").Value = '2.275.93'
$ws.Range("E3").Value = '  +2.94%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '498.08'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.38%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.38'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.58%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.20%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.529'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.23%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0959'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.23%  '
$ws.Range("E10").Value = '  +2.17%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.332'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.64%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.71'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.94%  '
$ws.Range("D13").Value = '2.681.67'
$ws.Range("E13").Value = '  +2.74%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '22.40'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.44%  '
$ws.Range("D15").Value = '54.273.62'
$ws.Range("E15").Value = '  +1.37%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000130'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.51%  '
$ws.Range("D17").Value = '2.276.58'
$ws.Range("E17").Value = '  +1.89%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.17'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +5.37%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.12'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.51%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '305.35'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.16%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.46'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.71%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '61.93'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.69%  '
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("D25").Value = '2.379.98'
$ws.Range("E25").Value = '  +1.73%  '
$ws.Range("E26").Value = '  +2.44%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.30'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.75%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '172.09'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.64%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.62'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.42%  '
$ws.Range("D30").Value = '0.0₃0688'
$ws.Range("E30").Value = '  +2.83%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.94'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.68%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.09'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.02%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '17.77'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.76%  '
$ws.Range("E35").Value = '  +0.30%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.926'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +10.99%  '
$ws.Range("E37").Value = '  +2.54%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.73'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.56%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '35.79'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.71%  '
$ws.Range("E40").Value = '  +1.19%  '
$ws.Range("E41").Value = '  +2.79%  '
$ws.Range("E42").Value = '  +2.78%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.99'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.33%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '126.52'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.86%  '
$ws.Range("E45").Value = '  +1.87%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0491'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.00%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.548'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.19%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '239.44'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.98%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.372'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.21%  '
$ws.Range("E50").Value = '  +2.64%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '10.77'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.03%  '
